$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I (I0) and J (IF). Copy the format from
# the existing H1 header (bold / centered / bordered) so the new headers
# reuse the same cell style as the rest of the header row, then set text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill data rows 2-27: column I is always 1, column J mirrors column H.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ws.Cells.Item($r, 8).Value()
}
